$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("A1:G1")
$rng.Merge()
$rng.Style = "Heading 2"
$rng.HorizontalAlignment = -4108
